$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add a new worksheet after the existing "Data" sheet
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Data2"

$ws2.Range("A1").Value = "Name "
$ws2.Range("B1").Value = "Address "
$ws2.Range("C1").Value = "Contact"
$ws2.Range("D1").Value = "City"
$ws2.Range("E1").Value = "PIN Code"

$ws2.Activate() | Out-Null
$ws2.Range("E1").Select() | Out-Null
